# Update workbook to reflect data refresh through 2021-10-15
# (commit: "Add data for 2021-10-23")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet (tab name) from "Through 2021-10-14" -> "Through 2021-10-15"
$ws.Name = "Through 2021-10-15"

# Update the header label in column B (shared string) to match the new "through" date
$ws.Range("B1").Value = "October 2021 (through October 15)"

# --- Updated / newly-added monthly carjacking counts by neighborhood ---

# Row 2 - Garfield Park
$ws.Range("L2").Value = 9    # October 2020: 8 -> 9
$ws.Range("AF2").Value = 4   # October 2018: 2 -> 4
$ws.Range("AZ2").Value = 4   # October 2016: 3 -> 4

# Row 3 - Austin
$ws.Range("B3").Value = 10   # October 2021 (through Oct 15): 7 -> 10
$ws.Range("L3").Value = 9    # October 2020: 7 -> 9
$ws.Range("AP3").Value = 4   # October 2017: 2 -> 4

# Row 4 - North Lawndale
$ws.Range("L4").Value = 7    # October 2020: 6 -> 7

# Row 5 - Chatham
$ws.Range("L5").Value = 1    # October 2020: new value

# Row 9 - Grand Crossing
$ws.Range("B9").Value = 4    # October 2021 (through Oct 15): 3 -> 4
$ws.Range("AZ9").Value = 3   # October 2016: 2 -> 3

# Row 10 - Roseland
$ws.Range("AP10").Value = 2  # October 2017: 1 -> 2

# Row 13 - South Shore
$ws.Range("AP13").Value = 5  # October 2017: 2 -> 5

# Row 14 - United Center
$ws.Range("L14").Value = 2   # October 2020: 1 -> 2

# Row 22 - Grand Boulevard
$ws.Range("AZ22").Value = 1  # October 2016: new value

# Row 24 - Chinatown
$ws.Range("B24").Value = 2   # October 2021 (through Oct 15): 1 -> 2

# Row 25 - River North
$ws.Range("B25").Value = 2   # October 2021 (through Oct 15): 1 -> 2

# Row 33 - Englewood
$ws.Range("B33").Value = 4   # October 2021 (through Oct 15): 3 -> 4

# Row 45 - West Lawn
$ws.Range("L45").Value = 1   # October 2020: new value

# Row 48 - Bridgeport
$ws.Range("L48").Value = 1   # October 2020: new value

# Row 57 - New City
$ws.Range("B57").Value = 2   # October 2021 (through Oct 15): 1 -> 2

# Row 59 - Irving Park
$ws.Range("L59").Value = 2   # October 2020: 1 -> 2

# Row 62 - Beverly
$ws.Range("L62").Value = 1   # October 2020: new value

# Row 64 - Brighton Park
$ws.Range("AZ64").Value = 1  # October 2016: new value
